$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "Última actualización: 10:36:50"
$ws.Range("A3").Value = "Total filas: 113"
$ws.Range("A62").Value = "07:49:32"
$ws.Range("C62").Value = "14_ABASTO"
$ws.Range("D62").Value = 88
$ws.Range("A63").Value = "08:38:24"
$ws.Range("C63").Value = "27_EL RETIRO"
$ws.Range("D63").Value = 39
$ws.Range("A88").Value = "10:36:50"
$ws.Range("B88").Value = "10:43"
$ws.Range("D88").Value = 7
$ws.Range("A89").Value = "08:45:31"
$ws.Range("B89").Value = "10:44"
$ws.Range("C89").Value = "11X44_ETCHEVERRY"
$ws.Range("D89").Value = 119
$ws.Range("A90").Value = "08:52:40"
$ws.Range("B90").Value = "10:46"
$ws.Range("C90").Value = "15_P INDUSTRIAL"
$ws.Range("D90").Value = 114
$ws.Range("A91").Value = "09:22:34"
$ws.Range("B91").Value = "10:53"
$ws.Range("D91").Value = 91
$ws.Range("A92").Value = "10:36:50"
$ws.Range("B92").Value = "10:55"
$ws.Range("C92").Value = "16_SANTA ANA"
$ws.Range("D92").Value = 19
$ws.Range("B93").Value = "10:56"
$ws.Range("C93").Value = "27_EL RETIRO"
$ws.Range("D93").Value = 52
$ws.Range("B94").Value = "10:57"
$ws.Range("C94").Value = "10_OLMOS"
$ws.Range("D94").Value = 95
$ws.Range("B95").Value = "10:59"
$ws.Range("C95").Value = "10_OLMOS"
$ws.Range("D95").Value = 55
$ws.Range("B96").Value = "11:01"
$ws.Range("C96").Value = "81_EL PELIGRO"
$ws.Range("D96").Value = 99
$ws.Range("A97").Value = "10:04:30"
$ws.Range("B97").Value = "11:03"
$ws.Range("C97").Value = "23_HERNANDEZ"
$ws.Range("D97").Value = 59
$ws.Range("A98").Value = "10:36:50"
$ws.Range("B98").Value = "11:06"
$ws.Range("C98").Value = "23_HERNANDEZ"
$ws.Range("D98").Value = 30
$ws.Range("A99").Value = "09:22:34"
$ws.Range("B99").Value = "11:10"
$ws.Range("C99").Value = "16_P MOR-SANTA ANA"
$ws.Range("D99").Value = 108
$ws.Range("A100").Value = "09:22:34"
$ws.Range("B100").Value = "11:14"
$ws.Range("C100").Value = "14_ABASTO"
$ws.Range("D100").Value = 112
$ws.Range("A101").Value = "09:22:34"
$ws.Range("B101").Value = "11:15"
$ws.Range("C101").Value = "15X38_ABASTO"
$ws.Range("D101").Value = 113
$ws.Range("A102").Value = "10:36:50"
$ws.Range("B102").Value = "11:25"
$ws.Range("C102").Value = "16_SANTA ANA"
$ws.Range("D102").Value = 49
$ws.Range("B103").Value = "11:29"
$ws.Range("C103").Value = "10_OLMOS"
$ws.Range("D103").Value = 85
$ws.Range("B104").Value = "11:29"
$ws.Range("C104").Value = "16_SANTA ANA"
$ws.Range("D104").Value = 85
$ws.Range("A105").Value = "10:36:50"
$ws.Range("B105").Value = "11:30"
$ws.Range("C105").Value = "215C_EL PATO"
$ws.Range("D105").Value = 54
$ws.Range("A106").Value = "10:04:30"
$ws.Range("B106").Value = "11:31"
$ws.Range("C106").Value = "215C_EL PATO"
$ws.Range("D106").Value = 87
$ws.Range("E106").Value = "LP1912"
$ws.Range("A107").Value = "10:04:30"
$ws.Range("B107").Value = "11:41"
$ws.Range("C107").Value = "215B_EL PATO"
$ws.Range("D107").Value = 97
$ws.Range("E107").Value = "LP1912"
$ws.Range("A108").Value = "10:04:30"
$ws.Range("B108").Value = "11:45"
$ws.Range("C108").Value = "15X38_ABASTO"
$ws.Range("D108").Value = 101
$ws.Range("E108").Value = "LP1912"
$ws.Range("A109").Value = "10:36:50"
$ws.Range("B109").Value = "11:48"
$ws.Range("C109").Value = "23_HERNANDEZ"
$ws.Range("D109").Value = 72
$ws.Range("E109").Value = "LP1912"
$ws.Range("A110").Value = "10:36:50"
$ws.Range("B110").Value = "11:52"
$ws.Range("C110").Value = "225_GOMEZ"
$ws.Range("D110").Value = 76
$ws.Range("E110").Value = "LP1912"
$ws.Range("A111").Value = "10:04:30"
$ws.Range("B111").Value = "11:53"
$ws.Range("C111").Value = "225_GOMEZ"
$ws.Range("D111").Value = 109
$ws.Range("E111").Value = "LP1912"
$ws.Range("A112").Value = "10:04:30"
$ws.Range("B112").Value = "11:58"
$ws.Range("C112").Value = "17_ROMERO"
$ws.Range("D112").Value = 114
$ws.Range("E112").Value = "LP1912"
$ws.Range("A113").Value = "10:36:50"
$ws.Range("B113").Value = "12:05"
$ws.Range("C113").Value = "11_ETCHEVERRY"
$ws.Range("D113").Value = 89
$ws.Range("E113").Value = "LP1912"
$ws.Range("A114").Value = "10:36:50"
$ws.Range("B114").Value = "12:10"
$ws.Range("C114").Value = "15_ABASTO"
$ws.Range("D114").Value = 94
$ws.Range("E114").Value = "LP1912"
$ws.Range("A115").Value = "10:36:50"
$ws.Range("B115").Value = "12:10"
$ws.Range("C115").Value = "16_P MOR-SANTA ANA"
$ws.Range("D115").Value = 94
$ws.Range("E115").Value = "LP1912"
$ws.Range("A116").Value = "10:36:50"
$ws.Range("B116").Value = "12:21"
$ws.Range("C116").Value = "215C_EL PATO"
$ws.Range("D116").Value = 105
$ws.Range("E116").Value = "LP1912"
$ws.Range("A117").Value = "10:36:50"
$ws.Range("B117").Value = "12:32"
$ws.Range("C117").Value = "14_ABASTO"
$ws.Range("D117").Value = 116
$ws.Range("E117").Value = "LP1912"
$ws.Range("A118").Value = "10:36:50"
$ws.Range("B118").Value = "12:34"
$ws.Range("C118").Value = "15_ABASTO"
$ws.Range("D118").Value = 118
$ws.Range("E118").Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = "Última actualización: 10:36:50"
$ws.Range("A3").Value = "Total filas: 22"
$ws.Range("A24").Value = "10:36:50"
$ws.Range("B24").Value = "11:30"
$ws.Range("D24").Value = 54
$ws.Range("B25").Value = "11:31"
$ws.Range("C25").Value = "215C_EL PATO"
$ws.Range("D25").Value = 87
$ws.Range("A26").Value = "10:04:30"
$ws.Range("B26").Value = "11:41"
$ws.Range("C26").Value = "215B_EL PATO"
$ws.Range("D26").Value = 97
$ws.Range("E26").Value = "LP1912"
$ws.Range("A27").Value = "10:36:50"
$ws.Range("B27").Value = "12:21"
$ws.Range("C27").Value = "215C_EL PATO"
$ws.Range("D27").Value = 105
$ws.Range("E27").Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "Última actualización: 10:36:50"
$ws.Range("A3").Value = "Total filas: 19"
$ws.Range("A23").Value = "10:36:50"
$ws.Range("B23").Value = "11:25"
$ws.Range("D23").Value = 49
$ws.Range("A24").Value = "10:04:30"
$ws.Range("B24").Value = "11:26"
$ws.Range("C24").Value = "215C_LA PLATA"
$ws.Range("D24").Value = 82
$ws.Range("E24").Value = "L6203"
